$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header: "intervention_type"; copy formatting from the
# existing header row (A1:J1) so it shares the same bold/border/centered style.
$ws.Range("K1").Value = "intervention_type"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Per-row intervention_type values. A handful of rows (14, 17, 69, 97, 124)
# have no intervention type and are intentionally left blank.
$ws.Range("K2").Value = "BIOLOGICAL"
$ws.Range("K3").Value = "DRUG"
$ws.Range("K4").Value = "DRUG"
$ws.Range("K5").Value = "DRUG"
$ws.Range("K6").Value = "OTHER"
$ws.Range("K7").Value = "PROCEDURE"
$ws.Range("K8").Value = "DEVICE"
$ws.Range("K9").Value = "DRUG"
$ws.Range("K10").Value = "OTHER"
$ws.Range("K11").Value = "DRUG"
$ws.Range("K12").Value = "DRUG"
$ws.Range("K13").Value = "DRUG"
# K14: left blank (no intervention_type)
$ws.Range("K15").Value = "DEVICE"
$ws.Range("K16").Value = "DRUG"
# K17: left blank (no intervention_type)
$ws.Range("K18").Value = "DRUG"
$ws.Range("K19").Value = "DRUG"
$ws.Range("K20").Value = "OTHER"
$ws.Range("K21").Value = "DRUG"
$ws.Range("K22").Value = "DRUG"
$ws.Range("K23").Value = "OTHER"
$ws.Range("K24").Value = "DRUG"
$ws.Range("K25").Value = "OTHER"
$ws.Range("K26").Value = "OTHER"
$ws.Range("K27").Value = "PROCEDURE"
$ws.Range("K28").Value = "DEVICE"
$ws.Range("K29").Value = "DEVICE"
$ws.Range("K30").Value = "BEHAVIORAL"
$ws.Range("K31").Value = "OTHER"
$ws.Range("K32").Value = "PROCEDURE"
$ws.Range("K33").Value = "PROCEDURE"
$ws.Range("K34").Value = "DEVICE"
$ws.Range("K35").Value = "PROCEDURE"
$ws.Range("K36").Value = "DRUG"
$ws.Range("K37").Value = "DEVICE"
$ws.Range("K38").Value = "BIOLOGICAL"
$ws.Range("K39").Value = "DRUG"
$ws.Range("K40").Value = "DEVICE"
$ws.Range("K41").Value = "RADIATION"
$ws.Range("K42").Value = "OTHER"
$ws.Range("K43").Value = "DRUG"
$ws.Range("K44").Value = "OTHER"
$ws.Range("K45").Value = "OTHER"
$ws.Range("K46").Value = "DEVICE"
$ws.Range("K47").Value = "PROCEDURE"
$ws.Range("K48").Value = "OTHER"
$ws.Range("K49").Value = "PROCEDURE"
$ws.Range("K50").Value = "OTHER"
$ws.Range("K51").Value = "DRUG"
$ws.Range("K52").Value = "DRUG"
$ws.Range("K53").Value = "DRUG"
$ws.Range("K54").Value = "BIOLOGICAL"
$ws.Range("K55").Value = "DIAGNOSTIC_TEST"
$ws.Range("K56").Value = "OTHER"
$ws.Range("K57").Value = "DRUG"
$ws.Range("K58").Value = "OTHER"
$ws.Range("K59").Value = "DEVICE"
$ws.Range("K60").Value = "DEVICE"
$ws.Range("K61").Value = "DRUG"
$ws.Range("K62").Value = "DRUG"
$ws.Range("K63").Value = "PROCEDURE"
$ws.Range("K64").Value = "PROCEDURE"
$ws.Range("K65").Value = "DRUG"
$ws.Range("K66").Value = "DEVICE"
$ws.Range("K67").Value = "DIAGNOSTIC_TEST"
$ws.Range("K68").Value = "OTHER"
# K69: left blank (no intervention_type)
$ws.Range("K70").Value = "DRUG"
$ws.Range("K71").Value = "PROCEDURE"
$ws.Range("K72").Value = "BIOLOGICAL"
$ws.Range("K73").Value = "OTHER"
$ws.Range("K74").Value = "DRUG"
$ws.Range("K75").Value = "BEHAVIORAL"
$ws.Range("K76").Value = "DRUG"
$ws.Range("K77").Value = "DRUG"
$ws.Range("K78").Value = "PROCEDURE"
$ws.Range("K79").Value = "DRUG"
$ws.Range("K80").Value = "BIOLOGICAL"
$ws.Range("K81").Value = "DRUG"
$ws.Range("K82").Value = "PROCEDURE"
$ws.Range("K83").Value = "DRUG"
$ws.Range("K84").Value = "OTHER"
$ws.Range("K85").Value = "OTHER"
$ws.Range("K86").Value = "DIAGNOSTIC_TEST"
$ws.Range("K87").Value = "PROCEDURE"
$ws.Range("K88").Value = "DEVICE"
$ws.Range("K89").Value = "DEVICE"
$ws.Range("K90").Value = "PROCEDURE"
$ws.Range("K91").Value = "DRUG"
$ws.Range("K92").Value = "OTHER"
$ws.Range("K93").Value = "BIOLOGICAL"
$ws.Range("K94").Value = "OTHER"
$ws.Range("K95").Value = "DRUG"
$ws.Range("K96").Value = "DRUG"
# K97: left blank (no intervention_type)
$ws.Range("K98").Value = "DEVICE"
$ws.Range("K99").Value = "DEVICE"
$ws.Range("K100").Value = "PROCEDURE"
$ws.Range("K101").Value = "DRUG"
$ws.Range("K102").Value = "DRUG"
$ws.Range("K103").Value = "OTHER"
$ws.Range("K104").Value = "DRUG"
$ws.Range("K105").Value = "DRUG"
$ws.Range("K106").Value = "OTHER"
$ws.Range("K107").Value = "PROCEDURE"
$ws.Range("K108").Value = "BEHAVIORAL"
$ws.Range("K109").Value = "PROCEDURE"
$ws.Range("K110").Value = "PROCEDURE"
$ws.Range("K111").Value = "PROCEDURE"
$ws.Range("K112").Value = "DEVICE"
$ws.Range("K113").Value = "DEVICE"
$ws.Range("K114").Value = "DEVICE"
$ws.Range("K115").Value = "DEVICE"
$ws.Range("K116").Value = "BEHAVIORAL"
$ws.Range("K117").Value = "DEVICE"
$ws.Range("K118").Value = "DEVICE"
$ws.Range("K119").Value = "DEVICE"
$ws.Range("K120").Value = "DEVICE"
$ws.Range("K121").Value = "DEVICE"
$ws.Range("K122").Value = "OTHER"
$ws.Range("K123").Value = "DRUG"
# K124: left blank (no intervention_type)
